$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update email address + hyperlink (was rubenzhito94@gmail.com) ---
$ws.Range("B2").Hyperlinks.Delete()
$ws.Range("B2").Value = "luisreinoso.03@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:luisreinoso.03@gmail.com") | Out-Null
$ws.Range("B2").Style = "Hipervínculo"

# --- Row 3: new socio record ---
$ws.Range("A3").Value = "RUVERLI TENAZOA ONORBE"
$ws.Range("B3").Value = "rubenzhito94@gmail.com "
$ws.Range("C3").Value = 73789548
$ws.Range("D3").Value = "San Lorenzo"
$ws.Range("F3").Value = "socio"

$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:rubenzhito94@gmail.com ") | Out-Null
$ws.Range("B3").Style = "Hipervínculo"

# --- Selection moves from D3 to B3 ---
$ws.Range("B3").Select() | Out-Null
